# Brasileirao 2023 data fix-up:
#  - rounds 15 (rows 145 & 151): the home/away goal counts were entered as
#    placeholder zeros; the match hadn't actually been played, so blank them
#    out again (keep the existing red "pending" cell style).
#  - round 21 (rows 202-211) and round 22 (rows 213-221): results have come
#    in, so fill in the home (D) / away (F) goal counts using the regular
#    body-cell style used by the rest of the sheet.
#  - round 22 row 212 is still unplayed, so just mark its D/F score cells
#    with a distinct "not played yet" style (no value) instead of leaving
#    them absent altogether.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Round 15: matches that turned out not to have been played after all -
# clear the placeholder 0 scores but leave the highlighted style in place.
# ---------------------------------------------------------------------------
$ws.Range("D145").ClearContents()
$ws.Range("F145").ClearContents()
$ws.Range("D151").ClearContents()
$ws.Range("F151").ClearContents()

# ---------------------------------------------------------------------------
# New cell style for "scheduled, not played" score cells (row 212), matching
# the look of the rest of the data but with no alignment override - gives us
# a brand new entry in cellXfs (count 6 -> 7).
# ---------------------------------------------------------------------------
$pendingStyle = $wb.Styles.Add("Placar_Pendente")
$pendingStyle.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# Round 21 & 22 results. Home goals (D) / away goals (F) for each match,
# keyed by row number. Row 212 has no result yet.
# ---------------------------------------------------------------------------
$homeGoals = @{
    202 = 1; 203 = 2; 204 = 3; 205 = 2; 206 = 3; 207 = 2; 208 = 2; 209 = 1; 210 = 3; 211 = 0;
    213 = 1; 214 = 1; 215 = 0; 216 = 0; 217 = 0; 218 = 2; 219 = 2; 220 = 1; 221 = 1
}
$awayGoals = @{
    202 = 1; 203 = 1; 204 = 1; 205 = 2; 206 = 0; 207 = 0; 208 = 0; 209 = 0; 210 = 0; 211 = 0;
    213 = 2; 214 = 0; 215 = 0; 216 = 0; 217 = 0; 218 = 0; 219 = 0; 220 = 1; 221 = 1
}

foreach ($row in 202..221) {
    $dCell = $ws.Range("D$row")
    $fCell = $ws.Range("F$row")

    if ($row -eq 212) {
        # Still to be played - style only, no score.
        $dCell.Style = "Placar_Pendente"
        $fCell.Style = "Placar_Pendente"
        continue
    }

    # Copy the regular body-cell format (same style already used by column A
    # on this row) onto D/F before writing the value, so the new cells pick
    # up the plain (non-highlighted) style used throughout the sheet.
    $ws.Range("A$row").Copy()
    $dCell.PasteSpecial(-4122)
    $dCell.Value = $homeGoals[$row]

    $ws.Range("A$row").Copy()
    $fCell.PasteSpecial(-4122)
    $fCell.Value = $awayGoals[$row]
}

$excel.CutCopyMode = 0
